$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 27, shifting current rows 27/28 down to 28/29.
$ws.Rows.Item(27).Insert()

# Row 27 now holds the new weekly entry; fill it in (copy the shared
# descriptive columns from the row that used to occupy this spot, then
# overwrite the values that actually changed).
$ws.Range("A27").Value = 6
$ws.Range("B27").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C27").Value = "Metropolitana"
$ws.Range("D27").Value = 44610
$ws.Range("E27").Value = 13
$ws.Range("F27").Value = "Fruta"
$ws.Range("G27").Value = 100102
$ws.Range("H27").Value = "Cítricos"
$ws.Range("I27").Value = 100102006
$ws.Range("J27").Value = "Pomelo"
$ws.Range("K27").Value = "Start Ruby"
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 12
$ws.Range("N27").Value = 190000
$ws.Range("O27").Value = 190000
$ws.Range("P27").Value = 190000
$ws.Range("Q27").Value = "$/bins (350 kilos)"
$ws.Range("R27").Value = "Región Metropolitana"
$ws.Range("S27").Value = 543
$ws.Range("T27").Value = 350

$ws.Range("D27").NumberFormat = "YYYY-MM-DD HH:MM:SS"
